# Update "想去人数" (F column) figures for events on the "展览" and
# "全部类型" sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" --------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$sheet1Updates = @{
    2  = 174
    4  = 399
    6  = 5022
    8  = 575
    9  = 873
    13 = 543
    16 = 1648
    17 = 1431
    18 = 729
    21 = 266
    22 = 475
    23 = 119
    24 = 1040
    27 = 1978
    28 = 154
    31 = 211
    36 = 258
    37 = 568
    38 = 73
    39 = 28
    40 = 29
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    3  = 174
    5  = 399
    7  = 30
    8  = 5022
    10 = 575
    13 = 873
    19 = 543
    23 = 1648
    24 = 1431
    25 = 729
    28 = 266
    30 = 475
    31 = 119
    32 = 1040
    34 = 1978
    35 = 154
    38 = 211
    42 = 258
    43 = 568
    44 = 73
    45 = 28
    46 = 29
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

$wb.Save()
